$wb = $excel.ActiveWorkbook

# "Scripts" sheet (sheet1.xml) holds the training run log.
$ws = $wb.Worksheets.Item("Scripts")

# Row 6 was the previous run (simulid 9999) -> update to the new run id 2728.
$ws.Range("A6").Value = 2728

# Row 7 previously had no run id in column A; this run recorded id 9652.
$ws.Range("A7").Value = 9652

# Move the active selection to B8, matching where the user left off entering data.
[void]$ws.Range("B8").Select()
